# Add a new user/password pair row (row 23) to the GKUser data sheet,
# and move the active selection to the new cell B23 (was B22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "userpassword"
$ws.Range("B23").Value = "yuwy"

$ws.Range("B23").Select()
